$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2406947890818859
$ws.Range("C2").Value = 0.4665012406947891
$ws.Range("J2").Value = 0.01240694789081886
$ws.Range("P2").Value = 0.1885856079404467
$ws.Range("S2").Value = 0.09181141439205956
$ws.Range("B3").Value = 0.01058201058201058
$ws.Range("J3").Value = 0.04232804232804233
$ws.Range("P3").Value = 0.708994708994709
$ws.Range("S3").Value = 0.2380952380952381
$ws.Range("J4").Value = 0.03448275862068965
$ws.Range("P4").Value = 0.6666666666666666
$ws.Range("S4").Value = 0.2988505747126437
$ws.Range("B6").Value = 0.064
$ws.Range("D6").Value = 0.008
$ws.Range("F6").Value = 0.052
$ws.Range("J6").Value = 0.296
$ws.Range("O6").Value = 0.04
$ws.Range("Q6").Value = 0.16
$ws.Range("R6").Value = 0.036
$ws.Range("S6").Value = 0.344
$ws.Range("B7").Value = 0.1015625
$ws.Range("D7").Value = 0.04296875
$ws.Range("E7").Value = 0.00390625
$ws.Range("F7").Value = 0.03515625
$ws.Range("J7").Value = 0.13671875
$ws.Range("O7").Value = 0.0234375
$ws.Range("Q7").Value = 0.171875
$ws.Range("R7").Value = 0.08203125
$ws.Range("S7").Value = 0.40234375
$ws.Range("B8").Value = 0.08893709327548807
$ws.Range("D8").Value = 0.02819956616052061
$ws.Range("F8").Value = 0.07158351409978309
$ws.Range("J8").Value = 0.1323210412147506
$ws.Range("O8").Value = 0.02169197396963124
$ws.Range("Q8").Value = 0.1735357917570499
$ws.Range("R8").Value = 0.0455531453362256
$ws.Range("S8").Value = 0.438177874186551
$ws.Range("B9").Value = 0.1102661596958175
$ws.Range("D9").Value = 0.01901140684410646
$ws.Range("F9").Value = 0.04182509505703422
$ws.Range("J9").Value = 0.1140684410646388
$ws.Range("O9").Value = 0.02661596958174905
$ws.Range("Q9").Value = 0.2243346007604563
$ws.Range("R9").Value = 0.05703422053231939
$ws.Range("S9").Value = 0.4068441064638783
$ws.Range("B10").Value = 0.1181873851806491
$ws.Range("D10").Value = 0.03612982241273729
$ws.Range("E10").Value = 0.001224739742804654
$ws.Range("F10").Value = 0.05756276791181874
$ws.Range("J10").Value = 0.1298224127372933
$ws.Range("O10").Value = 0.02633190447030006
$ws.Range("Q10").Value = 0.2308634415186773
$ws.Range("R10").Value = 0.04715248009797918
$ws.Range("S10").Value = 0.3527250459277403
$ws.Range("F11").Value = 0.002604166666666667
$ws.Range("G11").Value = 0.1380208333333333
$ws.Range("J11").Value = 0.08854166666666667
$ws.Range("K11").Value = 0.1510416666666667
$ws.Range("L11").Value = 0.6041666666666666
$ws.Range("S11").Value = 0.015625
$ws.Range("G12").Value = 0.7204724409448819
$ws.Range("J12").Value = 0.1968503937007874
$ws.Range("K12").Value = 0.007874015748031496
$ws.Range("L12").Value = 0.02755905511811024
$ws.Range("S12").Value = 0.04724409448818898
$ws.Range("G13").Value = 0.5531914893617021
$ws.Range("J13").Value = 0.3617021276595745
$ws.Range("S13").Value = 0.0851063829787234
$ws.Range("F15").Value = 0.03254437869822485
$ws.Range("H15").Value = 0.106508875739645
$ws.Range("I15").Value = 0.07396449704142012
$ws.Range("J15").Value = 0.3846153846153846
$ws.Range("K15").Value = 0.0650887573964497
$ws.Range("M15").Value = 0.008875739644970414
$ws.Range("N15").Value = 0.002958579881656805
$ws.Range("O15").Value = 0.07692307692307693
$ws.Range("S15").Value = 0.2485207100591716
$ws.Range("F16").Value = 0.04166666666666666
$ws.Range("H16").Value = 0.1401515151515151
$ws.Range("I16").Value = 0.1022727272727273
$ws.Range("J16").Value = 0.3825757575757576
$ws.Range("K16").Value = 0.1212121212121212
$ws.Range("M16").Value = 0.04166666666666666
$ws.Range("O16").Value = 0.05303030303030303
$ws.Range("S16").Value = 0.1174242424242424
$ws.Range("F17").Value = 0.01848739495798319
$ws.Range("H17").Value = 0.1210084033613445
$ws.Range("I17").Value = 0.1042016806722689
$ws.Range("J17").Value = 0.4487394957983193
$ws.Range("K17").Value = 0.1260504201680672
$ws.Range("M17").Value = 0.008403361344537815
$ws.Range("N17").Value = 0.001680672268907563
$ws.Range("O17").Value = 0.0453781512605042
$ws.Range("S17").Value = 0.1260504201680672
$ws.Range("F18").Value = 0.01388888888888889
$ws.Range("H18").Value = 0.125
$ws.Range("I18").Value = 0.1041666666666667
$ws.Range("J18").Value = 0.4236111111111111
$ws.Range("K18").Value = 0.09722222222222222
$ws.Range("M18").Value = 0.02777777777777778
$ws.Range("N18").Value = 0.006944444444444444
$ws.Range("O18").Value = 0.09722222222222222
$ws.Range("S18").Value = 0.1041666666666667
$ws.Range("F19").Value = 0.01941747572815534
$ws.Range("H19").Value = 0.1883495145631068
$ws.Range("I19").Value = 0.08543689320388349
$ws.Range("J19").Value = 0.370873786407767
$ws.Range("K19").Value = 0.1093851132686084
$ws.Range("M19").Value = 0.0168284789644013
$ws.Range("O19").Value = 0.09514563106796116
$ws.Range("S19").Value = 0.1145631067961165
